# Weekly data refresh: a new price observation (week of 2021-11-something,
# serial date 44518) is published and inserted into the data block at
# row 157, pushing the existing rows 157-244 down to 158-245.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 157 - this shifts rows 157:244 down
# to 158:245 (and any formatting on row 157 carries down with them).
$ws.Rows.Item(157).Insert()

# Populate the newly inserted row 157 with the new observation.
$ws.Cells.Item(157, 1).Value  = 9
$ws.Cells.Item(157, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(157, 3).Value  = "Metropolitana"
$ws.Cells.Item(157, 4).Value  = 44518
$ws.Cells.Item(157, 5).Value  = 13
$ws.Cells.Item(157, 6).Value  = 100112044
$ws.Cells.Item(157, 7).Value  = "Perejil"
$ws.Cells.Item(157, 8).Value  = "Sin especificar"
$ws.Cells.Item(157, 9).Value  = "Primera"
$ws.Cells.Item(157, 10).Value = 106
$ws.Cells.Item(157, 11).Value = 10000
$ws.Cells.Item(157, 12).Value = 12000
$ws.Cells.Item(157, 13).Value = 11000
$ws.Cells.Item(157, 14).Value = "`$/docena de atados"
$ws.Cells.Item(157, 15).Value = "Región Metropolitana"
$ws.Cells.Item(157, 16).Value = 3667
$ws.Cells.Item(157, 17).Value = 3
$ws.Cells.Item(157, 18).Value = "Hortaliza"
